$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-12 01:40:05"

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
